# Append the 2025-03-22 price row to every "Solar_Prices" sheet.
# Each sheet has a Date/Price table in columns A:B ending at row 20;
# we add row 21 with date "2025-03-22" and the same price as the prior
# day (row 20), keeping both cells as plain text (matching the existing
# inline-string-like text cells already in the sheet).

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-22"

# Sheet name -> price text to put in column B of the new row.
$priceBySheet = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"              = "1.19"
    "Cell Topcon 183mm"         = "0.298"
    "Module Topcon 183mm"       = "0.1"
    "Silver Rear_side"          = "5,399"
    "Silver Busbar front-side"  = "8,083"
    "Silver finger front-side"  = "8,133"
    "USD_CNY"                   = "7.2717"
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($priceBySheet.ContainsKey($name)) {
        $price = $priceBySheet[$name]

        $dateCell = $ws.Range("A21")
        $priceCell = $ws.Range("B21")

        # Force text storage so date-looking / numeric-looking strings are
        # not auto-converted into a real date serial / number, then reset
        # the style back to Normal so no stray number-format style sticks
        # to the cell (keeping it identical to the surrounding cells).
        $dateCell.NumberFormat = "@"
        $dateCell.Value = $newDate
        $dateCell.Style = "Normal"

        $priceCell.NumberFormat = "@"
        $priceCell.Value = $price
        $priceCell.Style = "Normal"
    }
}
